$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 433, shifting existing rows 433:460 down to 434:461
$ws.Rows.Item(433).Insert()

# Populate the newly inserted row with the new weekly price entry
$ws.Cells.Item(433, 1).Value = 10
$ws.Cells.Item(433, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(433, 3).Value = "La Araucanía"
$ws.Cells.Item(433, 4).Value = 45041
$ws.Cells.Item(433, 5).Value = 9
$ws.Cells.Item(433, 6).Value = 100114013
$ws.Cells.Item(433, 7).Value = "Zanahoria"
$ws.Cells.Item(433, 8).Value = "Sin especificar"
$ws.Cells.Item(433, 9).Value = "Primera"
$ws.Cells.Item(433, 10).Value = 100
$ws.Cells.Item(433, 11).Value = 6000
$ws.Cells.Item(433, 12).Value = 6000
$ws.Cells.Item(433, 13).Value = 6000
$ws.Cells.Item(433, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(433, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(433, 16).Value = 240
$ws.Cells.Item(433, 17).Value = 25
$ws.Cells.Item(433, 18).Value = "Hortaliza"
